$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change reflects a localization "handback" report being generated:
# the two files that were "Ready for handoff" (4ebb4f34-...md and
# f08c2243-...md) are now "Handed back: in sync with en-US" for both the
# zh-cn and de-de locales, with their Latest Target File / Latest Handback
# File / Latest Handback DateTime columns populated.
# ---------------------------------------------------------------------------

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $handedBack
$wsOverview.Range("C3").Value = $handedBack
$wsOverview.Range("B4").Value = $handedBack
$wsOverview.Range("C4").Value = $handedBack

# --- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = $handedBack
$wsZh.Range("E3").Value = "4ebb4f34-fee2-4036-ac66-779545ca2a1a.md"
$wsZh.Range("F3").Value = "4ebb4f34-fee2-4036-ac66-779545ca2a1a.66049ba1281fc4f9c2ad97299cc08db995fc76ba.zh-cn.xlf"
$wsZh.Range("G3").Value = "2016-03-10 05:03:05"

$wsZh.Range("B4").Value = $handedBack
$wsZh.Range("E4").Value = "4ebb4f34-fee2-4036-ac66-779545ca2a1a.md"
$wsZh.Range("F4").Value = "4ebb4f34-fee2-4036-ac66-779545ca2a1a.66049ba1281fc4f9c2ad97299cc08db995fc76ba.zh-cn.xlf"
$wsZh.Range("G4").Value = "2016-03-10 05:03:05"

# --- de-de sheet ------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = $handedBack
$wsDe.Range("E3").Value = "4ebb4f34-fee2-4036-ac66-779545ca2a1a.md"
$wsDe.Range("F3").Value = "4ebb4f34-fee2-4036-ac66-779545ca2a1a.66049ba1281fc4f9c2ad97299cc08db995fc76ba.de-de.xlf"
$wsDe.Range("G3").Value = "2016-03-10 05:03:19"

$wsDe.Range("B4").Value = $handedBack
$wsDe.Range("E4").Value = "4ebb4f34-fee2-4036-ac66-779545ca2a1a.md"
$wsDe.Range("F4").Value = "4ebb4f34-fee2-4036-ac66-779545ca2a1a.66049ba1281fc4f9c2ad97299cc08db995fc76ba.de-de.xlf"
$wsDe.Range("G4").Value = "2016-03-10 05:03:19"
